$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "used" / "have on hand" quantities for a few BOM line items.
# (Column H = qty already have; column I = qty still needed = B - H)
$ws.Range("H21").Value = 1
$ws.Range("H22").Value = 1
$ws.Range("H23").Value = 2

# Line item price correction and quantity update
$ws.Range("E24").Value = 6.9
$ws.Range("E24").NumberFormat = "0.00"
$ws.Range("H24").Value = 1

# Restore the view / selection to the top of the list
$ws.Range("A25:K25").Select()
